$d = $word.ActiveDocument

# Locate the list-item paragraph whose entire content is the old answer
# "3200" and turn it into the new explanatory sentence, underlining the
# paragraph mark (matches the target markup's <w:pPr><w:rPr><w:u .../>).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "3200") {
        $r = $p.Range

        # Underline the whole paragraph (text + its mark) first; this is the
        # only reliable way this COM host exposes paragraph-mark-only
        # underline formatting (it stamps w:pPr/w:rPr/w:u at the same time).
        $r.Font.Underline = 1

        # Remove the old "3200" text via Delete() (not Text = ""), then
        # retype the new sentence through a fresh Range positioned at the
        # paragraph start. Unlike re-assigning .Text on the old (now
        # underlined) range, inserting through a brand-new collapsed Range
        # here does not inherit the run-level underline, so only the
        # paragraph mark stays underlined - exactly matching the target.
        $textOnly = $d.Range($r.Start, $r.End - 1)
        $textOnly.Delete()

        $p2 = $d.Paragraphs.Item($p.Index)
        $insertPt = $d.Range($p2.Range.Start, $p2.Range.Start)
        $insertPt.InsertBefore("Tiene que se mayor al factor de carga maximo 4.0")

        break
    }
}
